$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the acronym table (A2:B13) alphabetically by column A (Name),
# replacing the previous sort which was by column B (Acronym).
$sortRange = $ws.Range("A2:B13")
$key1 = $ws.Range("A2:A13")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($key1, 0, 1, 0, 0) | Out-Null

$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = -4142
$ws.Sort.MatchCase = $false
$ws.Sort.Orientation = 1
$ws.Sort.SortMethod = 1
$ws.Sort.Apply()

# Update the active selection to reflect where the user clicked afterwards.
$ws.Range("B24").Select() | Out-Null
